# Update the "想去人数" (number of people interested) counts for two
# exhibition rows on both the "展览" and "全部类型" sheets.
#   F3: 15 -> 16
#   F4: 3  -> 4

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 16
    $ws.Range("F4").Value = 4
}
